# Updates meanEMG legmaxROM values for Hjemme passive output.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values (legmaxROM columns 15/16 repeated)
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 (CON) values
$ws.Range("B2").Value = 7.26485333301558
$ws.Range("C2").ClearContents()
$ws.Range("D2").Value = 6.0430975014234463
$ws.Range("E2").ClearContents()

# Row 3 (STR) values
$ws.Range("B3").Value = 4.7734771590408371
$ws.Range("C3").Value = -8.6193128364008444
$ws.Range("D3").Value = 3.6418213723349879
$ws.Range("E3").Value = -10.505396392868107

# Update the selection to match the new data extent
$ws.Range("B1:E3").Select()
